$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Approver")

# The approver "Bingo" access-code value had a typo (# instead of @) and is
# corrected here to the proper email-style value "Bingo@12345" across all
# three rows that reference it.
$ws.Range("B2").Value = "Bingo@12345"
$ws.Range("B3").Value = "Bingo@12345"
$ws.Range("B4").Value = "Bingo@12345"

# Column B is widened slightly to fit the corrected value, matching the
# onscreen best-fit behaviour seen after the edit.
$ws.Columns.Item(2).ColumnWidth = 11.166666666666666

# Last cell touched/selected by the user was B4.
$ws.Range("B4").Select()
